# Combine Florenceville and Bristol prior to 2009 (#4)
# Remove the separate "Bristol", "Florenceville", and "Florenceville-Bristol"
# rows from the inconsistent-municipalities table, since their pre-2009 data
# has been combined elsewhere. The Excel table auto-resizes as rows are
# deleted from within it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows bottom-up so earlier row numbers stay valid as we go.
# Row 8: Florenceville-Bristol (2000-2008 data, now redundant)
$ws.Rows.Item(8).Delete()
# Row 7: Florenceville
$ws.Rows.Item(7).Delete()
# Row 4: Bristol
$ws.Rows.Item(4).Delete()
